$p = $ppt.ActivePresentation

$dollar = [char]36

# -----------------------------------------------------------------------
# Change 1: "Como p-valor > 0.5, concluímos que temos uma " ->
#           "Como p-valor > 0.05, concluímos que temos uma "
# -----------------------------------------------------------------------
$oldPValor = "Como p-valor > 0.5, concluímos que temos uma "
$newPValor = "Como p-valor > 0.05, concluímos que temos uma "

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $tr = $shape.TextFrame.TextRange
            $runCount = $tr.Runs().Count
            for ($ri = 1; $ri -le $runCount; $ri++) {
                $run = $tr.Runs($ri)
                if ($run.Text -eq $oldPValor) {
                    $run.Text = $newPValor
                }
            }
        }
    }
}

# -----------------------------------------------------------------------
# Change 2: merge the split "R" / "$180.480,00" runs into a single run
#           reading "R$180.480,00" (keeping the bold+underline formatting
#           of the original "R" run) and leave the trailing "." run as-is.
# -----------------------------------------------------------------------
$restText = $dollar + "180.480,00"
$target = "R" + $restText
$restLen = $restText.Length

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $tr = $shape.TextFrame.TextRange
            $full = $tr.Text
            $idx = $full.IndexOf($target)
            if ($idx -ge 0) {
                $rStart = $idx + 1  # 1-based for Characters()

                # Everything after the leading "R" (i.e. "$180.480,00")
                # used to be its own run; empty it out so only the first
                # run (now holding "R") remains, then rewrite that run's
                # text to include the full amount.
                $rest = $tr.Characters($rStart + 1, $restLen)
                if ($rest.Text -eq $restText) {
                    $rest.Text = ""
                }

                $rRun = $tr.Characters($rStart, 1)
                if ($rRun.Text -eq "R") {
                    $rRun.Text = $target
                }
            }
        }
    }
}
